# Update the three input cells on Sheet1 (Gross Expenditures, Total M, and
# Total Labor Cost). Everything else on the sheet (percentages, technician
# hours, suggested standard added cost) is formula-driven off these three
# values and recalculates automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 12827.43
$ws.Range("D4").Value = 1575.06
$ws.Range("D5").Value = 5818.73
